$d = $word.ActiveDocument

# 1. Move the _GoBack bookmark: delete it from its original location
#    (between "Грюффел" and "о" in the title), and merge the title run.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
$d.Content.Find.Execute("Грюффело", $true, $false, $false, $false, $false, $true, 1, $false, "Грюффело", 1) | Out-Null

# 2. Spelling / punctuation corrections throughout the poem text.
$d.Content.Find.Execute("особенно если время обедат ей.", $true, $false, $false, $false, $false, $true, 1, $false, "особенно если время обедать  ей.", 1) | Out-Null
$d.Content.Find.Execute("Зайди ко мне в норку, перекусим чуть чуть!", $true, $false, $false, $false, $false, $true, 1, $false, "Зайди ко мне в норку, перекусим чуть-чуть!", 1) | Out-Null
$d.Content.Find.Execute("Как, ты не знаеш, так вот он какой:", $true, $false, $false, $false, $false, $true, 1, $false, "Как, ты не знаешь, так вот он какой:", 1) | Out-Null
$d.Content.Find.Execute("С большущеми клыками,", $true, $false, $false, $false, $false, $true, 1, $false, "С большущими клыками,", 1) | Out-Null
$d.Content.Find.Execute(" длинющеми когтями,", $true, $false, $false, $false, $false, $true, 1, $false, " длиннющими когтями,", 1) | Out-Null
$d.Content.Find.Execute("острейшеми зубами,", $true, $false, $false, $false, $false, $true, 1, $false, "острейшими зубами,", 1) | Out-Null
$d.Content.Find.Execute("иса, боишся ты зря!", $true, $false, $false, $false, $false, $true, 1, $false, "иса, боишься ты зря!", 1) | Out-Null
$d.Content.Find.Execute("Смотри ка, Мышонок, какой", $true, $false, $false, $false, $false, $true, 1, $false, "Смотри-ка, Мышонок, какой", 1) | Out-Null
$d.Content.Find.Execute("Я был бы рад как нибудь в другой раз,", $true, $false, $false, $false, $false, $true, 1, $false, "Я был бы рад как-нибудь в другой раз,", 1) | Out-Null
$d.Content.Find.Execute("ты не знаеш, так ты мне поверь:", $true, $false, $false, $false, $false, $true, 1, $false, "ты не знаешь, так ты мне поверь:", 1) | Out-Null
$d.Content.Find.Execute("остлявы колени, тяжолая лапа,", $true, $false, $false, $false, $false, $true, 1, $false, "остлявы колени, тяжeлая лапа,", 1) | Out-Null
$d.Content.Find.Execute("Мышоннок продолж", $true, $false, $false, $false, $false, $true, 1, $false, "Мышонок продолж", 1) | Out-Null
$d.Content.Find.Execute("едь Грюффало нет, его выдумал я!", $true, $false, $false, $false, $false, $true, 1, $false, "едь Грюффeло нет, его выдумал я!", 1) | Out-Null
$d.Content.Find.Execute("остлявы колени, тяжолая лапа,", $true, $false, $false, $false, $false, $true, 1, $false, "остлявы колени, тяжeлая лапа,", 1) | Out-Null
$d.Content.Find.Execute("утерброт с мышкой - как раз то что надо!", $true, $false, $false, $false, $false, $true, 1, $false, "утеpбpод с мышкой - как раз то что надо!", 1) | Out-Null
$d.Content.Find.Execute("не показалось кто то кричал!", $true, $false, $false, $false, $false, $true, 1, $false, "не показалось кто-то кричал!", 1) | Out-Null
$d.Content.Find.Execute("ам на тропинке, должен кто то гулять.", $true, $false, $false, $false, $false, $true, 1, $false, "ам на тропинке, должен кто-то гулять.", 1) | Out-Null

# 3. "Но та лиш" -> "Но та лишь", then re-insert the _GoBack bookmark
#    right after the newly added "ь" (its new home in the edited doc).
$r = $d.Content
$r.Find.Execute("Но та лиш", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$r.InsertAfter("ь")
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)

